$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new record as row 42
$ws.Range("A42").Value = "Jawhar"
$ws.Range("B42").Value = "Hafsa"
$ws.Range("E42").Value = "ww5A_WMAAAAJ"
$ws.Range("C42").Value = "Mohammed VI Polytechnic University"
$ws.Range("D42").Value = "Morocco"
$ws.Range("F42").Value = "M"

# Match the formatting used by the other rows in the Genre column (F)
$ws.Range("F2").Copy()
$ws.Range("F42").PasteSpecial(-4122)

$ws.Range("G42").Value = 1987
$ws.Range("H42").Value = "Médecine, Biologie et Sciences de la Santé"

# Update the visible selection / scroll position to match the saved view
$null = $ws.Activate()
$win = $excel.ActiveWindow
$null = $ws.Range("E39").Select()
$win.ScrollRow = 23
$win.ScrollColumn = 1
